$wb = $excel.ActiveWorkbook

# Update "想去人数" (attendance count) values in the "展览" and "全部类型" sheets
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 72
    $ws.Range("F4").Value = 45
}
